$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$b = $ws.Range("E8").Borders.Item(7)
$b.LineStyle = 1
$b.Weight = -4138
$b.ColorIndex = -4105
